$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: 10/05/2023 time-log entry for the Internship task.
$ws.Range("A10").Value = 45204
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("A10").HorizontalAlignment = -4131

$ws.Range("B10").Value = "Internship"

$ws.Range("C10").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Widen the Description column to fit the longer entry.
$ws.Columns.Item(3).ColumnWidth = 107.5
